$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B7").Value = "halosalsa"
$ws.Range("B8").Value = "halosalsa@gmail.com"

$ws.Range("B8").Select()
